# Update the "LS1-GA" (G) result column with refreshed run numbers for
# several instances. The dependent "%Over" formulas in column H (and the
# Mean/Stdev summary rows) recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 7542
$ws.Range("G4").Value = 935856
$ws.Range("G5").Value = 52643
$ws.Range("G7").Value = 109266
$ws.Range("G8").Value = 1651384
$ws.Range("G9").Value = 1400046
$ws.Range("G10").Value = 747333
$ws.Range("G11").Value = 842272
$ws.Range("G12").Value = 1208125
$ws.Range("G14").Value = 146393

$excel.CalculateFull()

$ws.Range("G20").Select() | Out-Null
